$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Описание товара"
$ws.Range("F1").Value = "Артикул"
$ws.Range("G1").Value = "Цена"
$ws.Range("H1").Value = "Наличие"
$ws.Range("I1").Value = "Размерность(шт, кг, тонна)"
